$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking values are not
# auto-converted to numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.819.97"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "1.550.70"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "205.69"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "21.48"
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "1.772.56"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "1.552.24"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").Value = "3.69"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "0.513"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "26.826.22"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "61.16"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "215.13"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "0.0₃0687"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "7.27"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "4.09"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "9.03"
$ws.Range("E23").Value = "  -4.21%  "
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").Value = "152.80"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").Value = "14.91"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "1.351.93"
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("E35").Value = "  -3.55%  "
$ws.Range("D36").Value = "2.27"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").Value = "0.924"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").Value = "0.0163"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("D39").Value = "0.524"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "5.63"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("D43").Value = "0.994"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "1.78"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "63.02"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").Value = "2.26"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").Value = "1.686.27"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "85.94"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("D51").Value = "0.0₇0973"
$ws.Range("E51").Value = "  -1.94%  "

# Restore default cell style (the NumberFormat tweak above is only a means
# to force text storage; the visible/serialized style should stay default).
$ws.Range("D2:D51").Style = "Normal"
